# Add two new event rows ("2021 Take Your Child to Work Day" and the
# structural-biology COVID-19 workshop) to the "events" sheet, matching
# the commit "Updated with TYCTW day link".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 - 2021 Take Your Child to Work Day
# Cells are written in column order F,G,H,I,J,K,B,A so that new shared
# strings land in the same order they were originally authored.
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "Virtual Event"
$ws.Range("G2").Value = "https://bioinformatics.niaid.nih.gov/take-your-child-to-work"
$ws.Range("H2").Value = "Virtual Event"
$ws.Range("I2").Value = "No"
$ws.Range("J2").Value = 'The NIH Office of Research Services hosts the annual Take Your Child to Work Day with an aim to inspire the next generation of NIH daughters and sons in grades 1-12 to explore career paths in science and public service at our nation’s medical research agency. 2021 is being held as a virtual event. Registration has closed, but resources provided by NIAID’s 3D Printing and Biovisualization Program are being made publicly available here.'
$ws.Range("K2").Value = "3D,visualization,STEM"
$ws.Range("B2").Value = "2021 Take Your Child to Work Day"
$ws.Range("A2").Value = "NIH ALL STAFF LISTSERV"

$ws.Range("C2").Value = 44308
$ws.Range("D2").Value = 44308
$ws.Range("E2").Value = 1619096400

# ---------------------------------------------------------------------
# Row 3 - Advances in COVID-19 Prevention and Treatment Enabled by
# Structural Biology Research
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "STRUCTBIOLIG LISTSERV"
$ws.Range("B3").Value = "Advances in COVID-19 Prevention and Treatment Enabled by Structural Biology Research"
$ws.Range("F3").Value = "Virtual Workshop"
$ws.Range("G3").Value = "https://www.aps.anl.gov/sites/www.aps.anl.gov/files/APS-Uploads/WK9%20Agenda.pdf"
$ws.Range("H3").Value = "Virtual Event"
$ws.Range("I3").Value = "No"
$ws.Range("J3").Value = "Broadly, the workshop will present areas where structural biology research, including macromolecular crystallography and cryoelectron microscopy, intersects with in vivo, in vitro, and in silico studies of SARS-CoV-2 and COVID-19. More precisely, the topics will include (a) viral biology, (b) vaccine, therapeutic, and diagnostic antibody studies, and (c) small-molecule drug discovery as it relates to viral proteases and other viral proteins. In addition, as this year's events emphasize the need for a coordinated, long-term strategy to prevent future pandemics of zoonotic origin, a broader One Health perspective on viral pathogens will be presented."
$ws.Range("K3").Value = "structural biology,crystallography,SARS-CoV-2,drug discovery"

$ws.Range("C3").Value = 44327
$ws.Range("D3").Value = 44328
$ws.Range("E3").Value = 1620741600

# ---------------------------------------------------------------------
# Row heights (wrapped description text makes these rows taller)
# ---------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 68
$ws.Rows.Item(3).RowHeight = 102

# ---------------------------------------------------------------------
# Column H ("eventType") now has data, so it gets an explicit best-fit
# width like the other text columns.
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 10.83

# ---------------------------------------------------------------------
# View state: scrolled right to column G, with J9 as the active cell.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$ws.Range("J9").Select()

Write-Host "done"
